$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: merge the two runs that make up the "TypeScript ... grandes" + "."
#         sentence into a single run (same visible text, one <w:r>).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$found1 = $r1.Find.Execute("TypeScript es un lenguaje de programación de código abierto con herramientas de programación orientada a objetos, muy favorable si se tienen proyectos grandes", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $afterGrandes = $r1.End
    $dotRng = $d.Range($afterGrandes, $afterGrandes + 1)
    if ($dotRng.Text -eq ".") {
        $dotRng.Delete()
        $insPoint1 = $d.Range($afterGrandes, $afterGrandes)
        $insPoint1.InsertAfter(".")
    }
}

# ---------------------------------------------------------------------------
# Edit 2: extend the "Básicamente es el intérprete ... de este." paragraph
#         with ", además administra paquetes o dependencias para ser
#         instalados." (as a second, distinct run) and move the hidden
#         "_GoBack" bookmark onto this paragraph.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("Básicamente es el intérprete de código de JavaScript y es el que permite la ejecución de este", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Work out which paragraph (by stable 1-based index) holds the match
    # *before* editing, while range positions are still trustworthy.
    $paraIndex = 0
    $i = 1
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Start -le $r2.Start -and $p.Range.End -ge $r2.End) {
            $paraIndex = $i
        }
        $i = $i + 1
    }

    $afterEste = $r2.End

    # Drop the trailing ". " (period + space) that used to close the sentence.
    $trailRng = $d.Range($afterEste, $afterEste + 2)
    if ($trailRng.Text -eq ". ") {
        $trailRng.Delete()
    }

    # Insert the new clause as its own run: toggling a character attribute
    # on/off forces a run boundary so it is not silently re-merged with the
    # preceding run even though the resulting formatting is identical.
    $insPoint2 = $d.Range($afterEste, $afterEste)
    $insPoint2.InsertAfter(", además administra paquetes o dependencias para ser instalados.")
    $newRunRng = $d.Range($afterEste, $afterEste + 66)
    $newRunRng.Font.Bold = 1
    $newRunRng.Font.Bold = 0

    # Move the (hidden) "_GoBack" bookmark so it wraps this whole paragraph
    # (re-fetched by stable index now that the paragraph has grown).
    if ($paraIndex -gt 0) {
        $paraNow = $d.Paragraphs.Item($paraIndex)
        $bmRange = $d.Range($paraNow.Range.Start, $paraNow.Range.End)
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}

$word.ActiveDocument.Saved = $false
